$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new "Created Room ID" value for column L (style s="8" cells)
$rows = @(2,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76)
$values = @("135227283","135233922","135234813","135235529","135236586","135237682","135238991","135239782","135240579","135241876","135242777","135243635","135244516","135245687","135246537","135247676","135249272","135250107","135251274","135252142","135253051","135254413","135255344","135256237","135257121","135258058","135258970","135259854","135261182","135262135","135263459","135264778","135265726","135266672","135267878","135268863","135270320","135271508","135272485","135273742","135274796","135275834","135277942","135279211","135280290","135276902","135281352","135282756","135284005","135285109","135286633","135287719","135289199","135312849","135315348","135317647","135320828","135323641","135326620","135328886","135331458","135291154","135303727","135305006","135292956","135296970","135300821","135308093","135309839")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $cell = $ws.Range("L" + $rows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
    $cell.NumberFormat = "General"
}

# Row 77 uses a distinct style (s="30"); a direct NumberFormat round-trip on
# that cell collapses its style index into the common one, so populate a
# scratch cell as text, copy/paste-values onto L77 (preserves L77's own
# style untouched), then clear the scratch cell.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "135311114"
$ws.Range("Z1").NumberFormat = "General"
$ws.Range("Z1").Copy()
$ws.Range("L77").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
